$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the current (pre-edit) row 32 values; this data will become
# --- the newly inserted row 33 once row 32 is refreshed with a later date.
$oldRow32 = @()
for ($col = 1; $col -le 20; $col++) {
    $oldRow32 += ,$ws.Cells.Item(32, $col).Value2
}

# --- Insert a new blank row at position 33; this shifts old rows 33..60
# --- down to 34..61, keeping all their data/formatting intact.
$ws.Cells.Item(33, 1).EntireRow.Insert()

# --- Update row 32's date (Fecha) to the newer reading; all other fields
# --- of row 32 (price/origin/etc.) remain unchanged.
$ws.Cells.Item(32, 4).Value = 44435

# --- Populate the newly inserted row 33 with the data that row 32 used to
# --- hold (same price/origin record, but the older date).
$ws.Cells.Item(33, 1).Value = $oldRow32[0]
$ws.Cells.Item(33, 2).Value = $oldRow32[1]
$ws.Cells.Item(33, 3).Value = $oldRow32[2]
$ws.Cells.Item(33, 4).Value = $oldRow32[3]
$ws.Cells.Item(33, 5).Value = $oldRow32[4]
$ws.Cells.Item(33, 6).Value = $oldRow32[5]
$ws.Cells.Item(33, 7).Value = $oldRow32[6]
$ws.Cells.Item(33, 8).Value = $oldRow32[7]
$ws.Cells.Item(33, 9).Value = $oldRow32[8]
$ws.Cells.Item(33, 10).Value = $oldRow32[9]
$ws.Cells.Item(33, 11).Value = $oldRow32[10]
$ws.Cells.Item(33, 12).Value = $oldRow32[11]
$ws.Cells.Item(33, 13).Value = $oldRow32[12]
$ws.Cells.Item(33, 14).Value = $oldRow32[13]
$ws.Cells.Item(33, 15).Value = $oldRow32[14]
$ws.Cells.Item(33, 16).Value = $oldRow32[15]
$ws.Cells.Item(33, 17).Value = $oldRow32[16]
$ws.Cells.Item(33, 18).Value = $oldRow32[17]
$ws.Cells.Item(33, 19).Value = $oldRow32[18]
$ws.Cells.Item(33, 20).Value = $oldRow32[19]

# --- Append a brand-new record as the new last row (62).
$newRow = 62
$ws.Cells.Item($newRow, 1).Value = 11
$ws.Cells.Item($newRow, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($newRow, 3).Value = "Bíobío"
$ws.Cells.Item($newRow, 4).Value = 44432
$ws.Cells.Item($newRow, 5).Value = 8
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100108
$ws.Cells.Item($newRow, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($newRow, 9).Value = 100108002
$ws.Cells.Item($newRow, 10).Value = "Mango"
$ws.Cells.Item($newRow, 11).Value = "Sin especificar"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 200
$ws.Cells.Item($newRow, 14).Value = 8000
$ws.Cells.Item($newRow, 15).Value = 8500
$ws.Cells.Item($newRow, 16).Value = 8250
$ws.Cells.Item($newRow, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item($newRow, 18).Value = "Brasil"
$ws.Cells.Item($newRow, 19).Value = 2062
$ws.Cells.Item($newRow, 20).Value = 4

# --- Match the date-formatted number format used by the other Fecha cells
# --- in column D (row 33 already inherits it from the Insert(); the
# --- freshly appended row needs it applied explicitly).
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item(32, 4).NumberFormat
